# --- Logins sheet: add two more user rows -------------------------------
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A4").Value = 'taf'
$ws1.Range("B4").Value = 'pass123'
$ws1.Range("C4").Value = 'faculty'

$ws1.Range("A5").Value = 'mmohmand'
$ws1.Range("B5").Value = 'password123'
$ws1.Range("C5").Value = 'faculty'

# Approximate "best fit" column widths that were present in the target file
$ws1.Range("A1").ColumnWidth = 10.9
$ws1.Range("B1").ColumnWidth = 12.65
$ws1.Range("C1").ColumnWidth = 6.2

# --- New "Applicants" sheet ----------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = 'Applicants'

$headers = @('id','jobid','first','last','address','email','phone','degree','school','role','company','yearsexp','education','exp','skills','communication','comments','attachment')
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

function Set-TextValue($range, $text) {
    $range.NumberFormat = '@'
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2 - mohmand application #1
$ws2.Cells.Item(2, 1).Value  = 2
$ws2.Cells.Item(2, 2).Value  = 1
$ws2.Cells.Item(2, 3).Value  = 'muj'
$ws2.Cells.Item(2, 4).Value  = 'mohmand'
$ws2.Cells.Item(2, 5).Value  = '10 Sunny Street, Toronto'
$ws2.Cells.Item(2, 6).Value  = 'mmohmand@algomau.ca'
Set-TextValue $ws2.Cells.Item(2, 7) '1234567890'
$ws2.Cells.Item(2, 8).Value  = 'computer science'
$ws2.Cells.Item(2, 9).Value  = 'Algoma University'
$ws2.Cells.Item(2, 10).Value = 'ceo'
$ws2.Cells.Item(2, 11).Value = 'google'
$ws2.Cells.Item(2, 12).Value = 5
$ws2.Cells.Item(2, 13).Value = 4
$ws2.Cells.Item(2, 14).Value = 3
$ws2.Cells.Item(2, 15).Value = 4
$ws2.Cells.Item(2, 16).Value = 5
$ws2.Cells.Item(2, 17).Value = "He's awesome!"
$ws2.Cells.Item(2, 18).Value = 'C:\Users\muj\Desktop\algoma\Y Daniel Liang - Introduction to Java Programming and Data Structures, Comprehensive Version-Pearson (2017).pdf'

# Row 3 - mohmand application #2
$ws2.Cells.Item(3, 1).Value  = 4
$ws2.Cells.Item(3, 2).Value  = 1
$ws2.Cells.Item(3, 3).Value  = 'muj'
$ws2.Cells.Item(3, 4).Value  = 'mohmand'
$ws2.Cells.Item(3, 5).Value  = '10 Sunny Street, Toronto'
$ws2.Cells.Item(3, 6).Value  = 'mmohmand@algomau.ca'
Set-TextValue $ws2.Cells.Item(3, 7) '1234567890'
$ws2.Cells.Item(3, 8).Value  = 'computer science'
$ws2.Cells.Item(3, 9).Value  = 'Algoma University'
$ws2.Cells.Item(3, 10).Value = 'ceo'
$ws2.Cells.Item(3, 11).Value = 'google'
$ws2.Cells.Item(3, 12).Value = 5
$ws2.Cells.Item(3, 13).Value = 4
$ws2.Cells.Item(3, 14).Value = 3
$ws2.Cells.Item(3, 15).Value = 4
$ws2.Cells.Item(3, 16).Value = 5
$ws2.Cells.Item(3, 17).Value = "He's awesome!"
$ws2.Cells.Item(3, 18).Value = 'C:\Users\muj\Desktop\algoma\Y Daniel Liang - Introduction to Java Programming and Data Structures, Comprehensive Version-Pearson (2017).pdf'

# Row 4 - jane smith application
$ws2.Cells.Item(4, 1).Value  = 5
$ws2.Cells.Item(4, 2).Value  = 1
$ws2.Cells.Item(4, 3).Value  = 'jane '
$ws2.Cells.Item(4, 4).Value  = 'smith'
$ws2.Cells.Item(4, 5).Value  = '123 fairway ave, sault ste marie, ontario'
$ws2.Cells.Item(4, 6).Value  = 'janesmith@google.com'
Set-TextValue $ws2.Cells.Item(4, 7) '555555555'
$ws2.Cells.Item(4, 8).Value  = 'computer science'
$ws2.Cells.Item(4, 9).Value  = 'Waterloo'
$ws2.Cells.Item(4, 10).Value = 'developer'
$ws2.Cells.Item(4, 11).Value = 'google'
$ws2.Cells.Item(4, 12).Value = 2
$ws2.Cells.Item(4, 13).Value = 5
$ws2.Cells.Item(4, 14).Value = 4
$ws2.Cells.Item(4, 15).Value = 3
$ws2.Cells.Item(4, 16).Value = 2
$ws2.Cells.Item(4, 17).Value = 'Good candiate. Recommended for interview with team.'
$ws2.Cells.Item(4, 18).Value = 'C:\Users\muj\Desktop\algoma\Y Daniel Liang - Introduction to Java Programming and Data Structures, Comprehensive Version-Pearson (2017).pdf'

[void]$ws2.Range("T9").Select()
